$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1, 2, 4, 1960513584, 7388667569),
    @(4, 2, 1, 7, 77755103, 738866764),
    @(4, 2, 1, 7, 211256629, 7388667768),
    @(4, 2, 1, 7, 29718549, 7388667812),
    @(4, 2, 1, 7, 225809045, 7388667836),
    @(4, 2, 1, 7, 209420883, 7388667847),
    @(4, 2, 1, 7, 211318748, 7388667974),
    @(4, 2, 1, 7, 244806218, 738866817),
    @(5, 3, 5, 12, 418257725, 7388668189)
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

$ws.Range("E11:F19").NumberFormat = "#,##0"

$ws.Range("A1:F19").Select()
